$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($old)
    $r.Text = $new
}

# 1. "solutions choisis" -> "solutions choisies" and add a trailing period
Replace-Text " (solutions choisis, documentation)" " (solutions choisies, documentation)."

# 2. Add trailing period after "... utilisation du matériel, …)"
Replace-Text "Manque d’explications techniques générales (tuto STM, utilisation du matériel, …)" `
             "Manque d’explications techniques générales (tuto STM, utilisation du matériel, …)."

# 3. Remove trailing space and add period "Pas d'inventaire et matériel mal rangé " (straight apostrophe, as in source doc)
Replace-Text "Pas d'inventaire et matériel mal rangé " "Pas d'inventaire et matériel mal rangé."

# 4. Add trailing period after "Pas de planning des tâches"
Replace-Text "Pas de planning des tâches" "Pas de planning des tâches."

# 5. Add trailing period after "(En décembre, certains choix techniques n’étaient pas encore faits)"
Replace-Text "(En décembre, certains choix techniques n’étaient pas encore faits)" `
             "(En décembre, certains choix techniques n’étaient pas encore faits)."

# 6. Add trailing period after "Procédures pour acheter du matériel"
Replace-Text "Procédures pour acheter du matériel" "Procédures pour acheter du matériel."

# 7. Fix typo "travaille" -> "travail"
Replace-Text "Présenter le travaille des années précédentes" "Présenter le travail des années précédentes"

# 8. Add trailing period after "Présenter la salle et les ressources (humaines et matérielles)"
Replace-Text "Présenter la salle et les ressources (humaines et matérielles)" `
             "Présenter la salle et les ressources (humaines et matérielles)."

# 9. Add trailing period after "... éviter la disqualification, etc. )"
Replace-Text "Etudier le sujet de la compétition. Relever les enjeux (se localiser, agir, communiquer, éviter la disqualification, etc. )" `
             "Etudier le sujet de la compétition. Relever les enjeux (se localiser, agir, communiquer, éviter la disqualification, etc. )."

# 10. Fix typo "coter" -> "côté"
Replace-Text "Ne pas mettre de coter la communication externe ET interne." `
             "Ne pas mettre de côté la communication externe ET interne."
